$wb = $excel.ActiveWorkbook

# --- Rename sheets (tab names encode updated timestamps) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16504778967632313"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778981282294"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16504778981292317"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16504778981772342"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778982532625"

# --- Sheet 1 (GNG): update stimulus file names ---
$ws1.Range("B2").Value = "go_stims-16504778967252347.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778967462292.csv"
$ws1.Range("B4").Value = "go_stims-165047789674823.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778967622664.csv"

# --- Sheet 2 (NB): update stimulus file names ---
$ws2.Range("B2").Value = "ZB-match_1-1650477896811229.csv"
$ws2.Range("B3").Value = "TB-165047789795223.csv"
$ws2.Range("B4").Value = "TB-16504778980212657.csv"
$ws2.Range("B5").Value = "ZB-match_4-1650477896773234.csv"
$ws2.Range("B6").Value = "OB-1650477897166232.csv"
$ws2.Range("B7").Value = "TB-16504778981072323.csv"
$ws2.Range("B8").Value = "ZB-match_0-16504778969502285.csv"
$ws2.Range("B9").Value = "OB-16504778974632626.csv"
$ws2.Range("B10").Value = "OB-16504778970052655.csv"

# --- Sheet 4 (TOL): update stimulus file names ---
$ws4.Range("B2").Value = "MM_stims-1650477898144264.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778981312313.csv"
$ws4.Range("B4").Value = "MM_stims-16504778981602643.csv"
$ws4.Range("B5").Value = "ZM_stims-1650477898144264.csv"
$ws4.Range("B6").Value = "MM_stims-1650477898176264.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778981612337.csv"

# --- Sheet 5 (vSAT): update stimulus file names ---
$ws5.Range("B2").Value = "vSAT_stims-16504778982222645.csv"
$ws5.Range("B3").Value = "SAT_stims-1650477898180231.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778982062647.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778982382667.csv"
